# First ready to deploy Network ARM template
# Rebuild Sheet1 to the new 3-column (A/B/C) layout described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Start clean: drop old merges + contents/formatting so the new layout
# (rows/columns shifted around) doesn't inherit stale merges or styles. ---
$ws.Cells.UnMerge()
$ws.Cells.Clear()

# --- Row 1: ResourceGroup header (A1:B1) + PSScript header (C1) ---
$ws.Range("A1").Value = "ResourceGroup"
$ws.Range("C1").Value = "PSScript"

# --- Row 2: Network1 / Network2 sub-headers + ARM Network label (merged C2:C4) ---
$ws.Range("A2").Value = "Network1"
$ws.Range("B2").Value = "Network2"
$ws.Range("C2").Value = "ARM Network"

# --- Row 3: Subnet1 / Subnet2 ---
$ws.Range("A3").Value = "Subnet1"
$ws.Range("B3").Value = "Subnet2"

# --- Row 4: SGRule (merged A4:B4) ---
$ws.Range("A4").Value = "SGRule"

# --- Row 5: LB1 / LB2 + ARM LB label (merged C5:C9) ---
$ws.Range("A5").Value = "LB1"
$ws.Range("B5").Value = "LB2"
$ws.Range("C5").Value = "ARM LB"

# --- Row 6: PublicIP1 / PublicIP2 ---
$ws.Range("A6").Value = "PublicIP1"
$ws.Range("B6").Value = "PublicIP2"

# --- Row 7: BEPool1 / BEPool2 ---
$ws.Range("A7").Value = "BEPool1"
$ws.Range("B7").Value = "BEPool2"

# --- Row 8: LBRule1 / LBRule2 ---
$ws.Range("A8").Value = "LBRule1"
$ws.Range("B8").Value = "LBRule2"

# --- Row 9: LB1NATRDP / LB2NATRDP ---
$ws.Range("A9").Value = "LB1NATRDP"
$ws.Range("B9").Value = "LB2NATRDP"

# --- Row 10: NIC1-2 / NIC3 + ARM VM label (merged C10:C12) ---
$ws.Range("A10").Value = "NIC1-2"
$ws.Range("B10").Value = "NIC3"
$ws.Range("C10").Value = "ARM VM"

# --- Row 11: VM1-2 / VM3 ---
$ws.Range("A11").Value = "VM1-2"
$ws.Range("B11").Value = "VM3"

# --- Row 12: AvSet ---
$ws.Range("A12").Value = "AvSet"

# ===================== Formatting ======================

# Give every used cell A1:C12 the thin box border first.
$ws.Range("A1:C12").Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$ws.Range("A1:C12").Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

# Headers A1:B1 and A4:B4 -> centered (merged).
$ws.Range("A1:B1").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$ws.Range("A4:B4").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

# Sub-header row A2:B2 -> left aligned.
$ws.Range("A2:B2").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft

# C2:C3 carry the scientific number format used by the template style.
$ws.Range("C2:C3").NumberFormat = "0.00E+00"
$ws.Range("C2:C3").WrapText = $false

# Remaining formatted-but-unstyled cells (touch alignment so the xf records
# the applyAlignment flag exactly like the authored workbook).
$ws.Range("C4:C12").WrapText = $false

# --- Merges ---
$ws.Range("A1:B1").Merge()
$ws.Range("A4:B4").Merge()
$ws.Range("C2:C4").Merge()
$ws.Range("C5:C9").Merge()
$ws.Range("C10:C12").Merge()

# --- Column C width (new column added alongside the template) ---
$ws.Columns("C").ColumnWidth = 11.6

# --- Selection matches the re-saved workbook state ---
$ws.Range("A4:B4").Select()
